# nexial-script.xlsx template update
#
# [web] openInTab(name,url): NEW command to open a url in another tab.
# Also registers check(name)/checkByLocator(locator) and
# uncheck(name)/uncheckByLocator(locator) as new `desktop` commands.
#
# The "#system" sheet holds, per-column, an alphabetically sorted list of
# command names that each feed a named range (used for data-validation
# dropdowns elsewhere in the workbook). Column H backs the "desktop" name,
# column AA backs the "web" name. We insert the new command names in their
# correct alphabetical slot, push the remainder of the column down, and
# then grow the named range to cover the extra rows.
#
# NOTE: Range.Insert() on this host shifts the *entire* row (every column),
# so to keep columns independent we shift values manually, one column at a
# time, from the bottom up. (Also: this host's PowerShell doesn't bind
# named `-param` args on custom functions, so positional args are used.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

function Shift-ColumnDown($col, $fromRow, $toRow, $by) {
    # col: 1-based column index (H = 8, AA = 27)
    # moves values in rows fromRow..toRow down by $by rows, bottom-up
    for ($r = $toRow; $r -ge $fromRow; $r--) {
        $v = $ws.Cells.Item($r, $col).Value()
        $ws.Cells.Item($r + $by, $col).Value = $v
    }
}

# --- column H ("desktop") ------------------------------------------------
# insert check(name) / checkByLocator(locator) before clear(locator) @ H26
Shift-ColumnDown 8 26 107 2
$ws.Cells.Item(26, 8).Value = "check(name)"
$ws.Cells.Item(27, 8).Value = "checkByLocator(locator)"

# insert uncheck(name) / uncheckByLocator(locator) before useApp(appId),
# which (after the shift above) now sits at H102
Shift-ColumnDown 8 102 109 2
$ws.Cells.Item(102, 8).Value = "uncheck(name)"
$ws.Cells.Item(103, 8).Value = "uncheckByLocator(locator)"

# --- column AA ("web") ---------------------------------------------------
# insert openInTab(name,url) before refresh() @ AA85
Shift-ColumnDown 27 85 152 1
$ws.Cells.Item(85, 27).Value = "openInTab(name,url)"

# --- grow the named ranges to cover the newly added rows ------------------
$wb.Names.Item("desktop").RefersTo = "='#system'!`$H`$2:`$H`$111"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$153"

Write-Output "desktop/web command lists updated"
